$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" value -> literal text "false" (not boolean FALSE).
# A plain .Value="false" auto-types to a Boolean in Excel, so write it as a
# formula yielding the text "false" then paste-special as values-only; this
# collapses the formula back to a literal shared-string cell with no style
# change (matches a typed-as-text "false").
$b7 = $ws.Range("B7")
$b7.Formula = '="false"'
$b7.Copy()
$b7.PasteSpecial(-4163)

$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
$ws.Range("B17").Value = "Cardiorespiratory fitness categories based on VO2max percentiles"
